# Add two new completed books to the "Completed" sheet: Cleopatra and African Samurai.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Row 126: Cleopatra
$ws.Cells.Item(126, 1).Value = "Cleopatra"
$ws.Cells.Item(126, 2).Value = "Stacy Schiff"
$ws.Cells.Item(126, 3).NumberFormat = "m/d/yy"
$ws.Cells.Item(126, 3).Value = Get-Date -Year 2020 -Month 9 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(126, 4).NumberFormat = "m/d/yy"
$ws.Cells.Item(126, 4).Value = Get-Date -Year 2020 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(126, 5).Value = "biography;egypt;rome;history;alexandria;cleopatra;julius caesar;mark antony;alexander the great;war;politics"
$ws.Cells.Item(126, 6).Value = "Audio"
$ws.Cells.Item(126, 7).Value = "14 Hours 17 Mins"
$ws.Cells.Item(126, 8).Value = 3

# Row 127: African Samurai
$ws.Cells.Item(127, 1).Value = "African Samurai"
$ws.Cells.Item(127, 2).Value = "Thomas Lockley"
$ws.Cells.Item(127, 3).NumberFormat = "m/d/yy"
$ws.Cells.Item(127, 3).Value = Get-Date -Year 2020 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(127, 4).NumberFormat = "m/d/yy"
$ws.Cells.Item(127, 4).Value = Get-Date -Year 2020 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(127, 5).Value = "biography;history;japan;samurai;slavery;war"
$ws.Cells.Item(127, 6).Value = "Audio"
$ws.Cells.Item(127, 7).Value = "10 Hours 13 Mins"
$ws.Cells.Item(127, 8).Value = 3

# Match the author's final view/selection state as closely as this host supports.
$ws.Activate()
$ws.Range("A128").Select()
